$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update "progress" values in column D (and a couple in column C) for the
#    new "ver-venta" (VerProductoController) / semana 4 rollout.
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = 100
$ws.Range("C3").Value = 70
$ws.Range("D3").Value = 80
$ws.Range("D4").Value = 88
$ws.Range("D5").Value = 60
$ws.Range("D6").Value = 60
$ws.Range("D7").Value = 100
$ws.Range("D8").Value = 100
$ws.Range("D9").Value = 85
$ws.Range("D10").Value = 85
$ws.Range("D11").Value = 96
$ws.Range("D12").Value = 80
$ws.Range("D13").Value = 90
$ws.Range("D14").Value = 100

# ---------------------------------------------------------------------------
# 2. Consolidate the duplicated conditional-formatting "traffic light" rules.
#    Originally the 3-rule "cellIs between" banding was duplicated per block
#    (B2:E13, then again separately for B14, C14, D14 and E14). Extend the
#    first banding rule to cover the whole table (B2:E14) and drop the
#    now-redundant per-cell copies for row 14.
# ---------------------------------------------------------------------------
$fcs = $ws.Cells.FormatConditions

function Get-CellIsGroupIndices($fcs, $addr) {
    $idx = @()
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fc = $fcs.Item($i)
        if ($fc.Type -eq 1 -and $fc.AppliesTo.Address() -eq $addr) {
            $idx += $i
        }
    }
    return $idx
}

# Remove the redundant "cellIs" banding rules left over on B14 / C14 / D14.
foreach ($addr in @('$B$14', '$C$14', '$D$14')) {
    $idx = Get-CellIsGroupIndices $fcs $addr
    for ($j = $idx.Count - 1; $j -ge 0; $j--) {
        $fcs.Item($idx[$j]).Delete()
    }
}

# Remove the old B2:E13-only banding rule (it will be replaced by extending
# the E14 banding rule below, which keeps the original dxf/priority order).
$idx = Get-CellIsGroupIndices $fcs '$B$2:$E$13'
for ($j = $idx.Count - 1; $j -ge 0; $j--) {
    $fcs.Item($idx[$j]).Delete()
}

# Extend the remaining (former E14-only) "cellIs" banding rule so that it
# now covers the entire data range B2:E14.
$idx = Get-CellIsGroupIndices $fcs '$E$14'
if ($idx.Count -gt 0) {
    $fcs.Item($idx[0]).ModifyAppliesToRange($ws.Range("B2:E14"))
}

# ---------------------------------------------------------------------------
# 3. Restore the cell selection to just below the table, where the user left
#    off after entering the new data.
# ---------------------------------------------------------------------------
$ws.Range("D15").Select() | Out-Null
